# Update "want to go" (想去人数) counts for two events that appear on
# both the "展览" sheet and the "全部类型" sheet.
#   F2: 3279 -> 3288
#   F4: 56   -> 57
#   F5: 1247 -> 1279

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 3288
    $ws.Range("F4").Value = 57
    $ws.Range("F5").Value = 1279
}
